# "Empty Level 6-7 added"
#
# The underlying edit swaps the level data (level name / objective / goal,
# columns C:E, including their cell formatting) that lives in row 5
# ("Level 6" - Мечта шейха) with the row 6 data ("Level 7" - Пустыня), so
# that the desert level moves up to row 5 and the "Мечта шейха" level moves
# down to row 6. Columns A/B (Build index / Level number) and F:I (tank
# info) are left untouched.
#
# We do this the way a user would in Excel: copy row 5's C:E block out to a
# scratch cell, copy row 6's C:E block into row 5, then drop the saved
# block into row 6 - which carries over both values and formatting
# (number format / alignment) together, matching the style swap seen in
# the diff (D5 goes from style 2 to style 3, D6 from style 3 to style 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("C200:E200")

$ws.Range("C5:E5").Copy($scratch)
$ws.Range("C6:E6").Copy($ws.Range("C5:E5"))
$scratch.Copy($ws.Range("C6:E6"))
$scratch.Clear()

# Match the saved selection/active cell from the edited workbook.
$ws.Range("D5").Select()
